# SVN Revision #7055 - Update Std design system maps and fan power per
# 3/3 CEC NACM system map document with corrections to fan power table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Table 2 (rows 14-18): system-map column headers / codes were
# re-ordered & relabeled per the updated CEC NACM system map, and a new
# system column (Sys9 / "9 - HEATVENT") was appended in column L.
# ---------------------------------------------------------------------

# Row 14 - descriptive headers (D:L)
$ws.Range("D14").Value = "3a – SZAC"
$ws.Range("E14").Value = "3b – SZHP (no furnace)"
$ws.Range("F14").Value = "3c – SZDFHP (with furnace)"
$ws.Range("G14").Value = "7a – SZVAVAC "
$ws.Range("H14").Value = "7b – SZVAVHP"
$ws.Range("I14").Value = "7c – SZVAVDFHP (with furnace)"
$ws.Range("J14").Value = "5 – PVAV"
$ws.Range("K14").Value = "6 – VAV"
$ws.Range("L14").Value = "9 – HEATVENT"

# row 14 header cells D:L had a distinct style (wrap/bold-ish) in the old
# layout; the refreshed table uses the plain default style instead (only
# A14 keeps its original style), and the row no longer needs the explicit
# 30pt row height - let it autofit back to the default.
$ws.Range("D14:L14").Style = "Normal"
$ws.Rows.Item(14).AutoFit()

# Row 15 - short system codes (D:L)
$ws.Range("D15").Value = "Sys3a"
$ws.Range("E15").Value = "Sys3b"
$ws.Range("F15").Value = "Sys3c"
$ws.Range("G15").Value = "Sys7a"
$ws.Range("H15").Value = "Sys7b"
$ws.Range("I15").Value = "Sys7c"
$ws.Range("J15").Value = "Sys5"
$ws.Range("K15").Value = "Sys6"
$ws.Range("L15").Value = "Sys9"

# Row 16 - "<=5000" fan power index row, now with corrected values and the
# new column L
$ws.Range("D16").Value = 0.802
$ws.Range("E16").Value = 0.744
$ws.Range("F16").Value = 0.802
$ws.Range("G16").Value = 0.802
$ws.Range("H16").Value = 0.744
$ws.Range("I16").Value = 0.802
$ws.Range("J16").Value = 1
$ws.Range("K16").Value = 0.977
$ws.Range("L16").Value = 0.616

# Row 17 - "<=10000" fan power index row
$ws.Range("D17").Value = 0.78
$ws.Range("E17").Value = 0.72
$ws.Range("F17").Value = 0.78
$ws.Range("G17").Value = 0.78
$ws.Range("H17").Value = 0.72
$ws.Range("I17").Value = 0.78
$ws.Range("J17").Value = 1.022
$ws.Range("K17").Value = 1.013
$ws.Range("L17").Value = 0.62

# Row 18 - ">10000" fan power index row
$ws.Range("D18").Value = 0.748
$ws.Range("E18").Value = 0.676
$ws.Range("F18").Value = 0.748
$ws.Range("G18").Value = 0.748
$ws.Range("H18").Value = 0.676
$ws.Range("I18").Value = 0.748
$ws.Range("J18").Value = 0.964
$ws.Range("K18").Value = 0.947
$ws.Range("L18").Value = 0.605

# ---------------------------------------------------------------------
# Window/selection cosmetics captured in the saved file
# ---------------------------------------------------------------------
$ws.Range("I20:I22").Select() | Out-Null
$excel.ActiveWindow.Zoom = 70

Write-Host "Fan power allowance tables updated."
